# DDAf_2022_Tableau_annexe_Tab28.xlsx - "Add files via upload" edit
#
# Fixes a handful of "Etats" -> "États" typos in the region-group labels,
# adds the missing "RDM, pays en développement sans littoral" label (row 96
# had been wrongly re-using the "Afrique, ..." label from row 95), tweaks a
# couple of figures, and rewords the disclaimer paragraph.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix accented "États" typos + add the missing "RDM, ..." label ---
$ws.Range("B93").Value = "Afrique, petits États insulaires en développement"
$ws.Range("B94").Value = "RDM, petits États insulaires en développement"

# Row 96 previously (incorrectly) duplicated row 95's "Afrique, pays en
# développement sans littoral" label; it should read "RDM, ...".
$ws.Range("B96").Value = "RDM, pays en développement sans littoral"

$ws.Range("B97").Value = "Afrique, États fragiles"
$ws.Range("B98").Value = "RDM, États fragiles"

# --- Small data correction ---
$ws.Range("C92").Value = 972

# --- Reworded disclaimer paragraph ---
$ws.Range("A104").Value = "Responsabilité : Ce tableau ainsi que toutes les données qu'il peut comprendre, sont sans préjudice du statut de tout territoire, de la souveraineté s'exerçant sur ce dernier, du tracé des frontières et limites internationales, et du nom de tout territoire, ville ou région."
